# Update the cached "date" field text (datetimeFigureOut placeholder) on the
# slide master and on every slide layout from 2020-04-08 -> 2020-04-09, and
# update the subtitle on slide 1 from "Lakshmi & Paddy" to
# "Lakshmi & Padmanabhan".

$p = $ppt.ActivePresentation

$oldDate = "2020-04-08"
$newDate = "2020-04-09"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePh = $true
            }
        } catch {
            $isDatePh = $false
        }
        if ($isDatePh) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master's "Date Placeholder" shape.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's "Date Placeholder" shape.
$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# Any slide that overrides the date placeholder itself (defensive; none do
# in this deck, but keeps the script correct if that ever changes).
for ($S = 1; $S -le $p.Slides.Count; $S++) {
    Update-DatePlaceholder $p.Slides.Item($S).Shapes
}

# Slide 1 subtitle text.
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $sh = $s1.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq "Lakshmi & Paddy") {
            $tr.Text = "Lakshmi & Padmanabhan"
        }
    }
}
